# Rewrite the pairwise-comparison labels from "with" to the abbreviated "w/"
# on both sheets (affects the shared strings used in B1 and D1 headers).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AUC_G_Log_Type_Increase")
$ws1.Range("B1").Value = "Naive | Selfreport w/ App"
$ws1.Range("D1").Value = "Selfreport w/ App | App"

$ws2 = $wb.Worksheets.Item("AUC_I_Log_Type_Increase")
$ws2.Range("B1").Value = "Naive | Selfreport w/ App"
$ws2.Range("D1").Value = "Selfreport w/ App | App"

# The workbook was resaved with "AUC_G_Log_Type_Increase" as the active
# sheet (selection reset to its default), instead of "AUC_I_Log_Type_Increase".
$ws2.Range("A1").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
